$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Ligand average / total expression values (new TPM values)
$ws.Range("G2").Value = 1.136164666666667
$ws.Range("H2").Value = 3.408494

# Row 2 - derived specificity values (recomputed after the TPM update)
$ws.Range("I2").Value = 0.01984988071065505
$ws.Range("J2").Value = 0.01984988071065505
$ws.Range("Q2").Value = 0.2834522546477778
$ws.Range("R2").Value = 2.55107029183
$ws.Range("S2").Value = 0.01984988071065505
$ws.Range("T2").Value = 0.01984988071065505

# Row 3 - derived specificity values (recomputed after the TPM update)
$ws.Range("I3").Value = 0.9624690187571424
$ws.Range("J3").Value = 0.9624690187571425
$ws.Range("S3").Value = 0.9624690187571424
$ws.Range("T3").Value = 0.9624690187571425

# Row 4 - derived specificity values (recomputed after the TPM update)
$ws.Range("I4").Value = 0.01768110053220253
$ws.Range("J4").Value = 0.01768110053220253
$ws.Range("S4").Value = 0.01768110053220253
$ws.Range("T4").Value = 0.01768110053220253
